$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header total count
$ws.Range("B1").Value = "Count (Total: 237)"

# Update the weekly triaged issue counts
$ws.Range("B2").Value = 143
$ws.Range("B3").Value = 80
$ws.Range("B4").Value = 14
